$d = $word.ActiveDocument

function Replace-Unique($find, $replace) {
    # Content-wide, case-sensitive, whole-document replace of a string that is
    # expected to occur exactly once in the document.
    $rng = $d.Content
    $null = $rng.Find.Execute($find, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $replace, 2)
}

# 1. City name in the dateline: Sobral -> Maceió
Replace-Unique "Sobral" "Maceió"

# 2. Requester block: two occurrences of the same company name become two
#    different strings (trade name, then legal/"razão social" name). They
#    must be resolved one at a time, in document order, because a blanket
#    Find/Replace-All would turn both into the same text.
$oldCompany = "DAQUIBRASIL SINAL DE FIBRA OTICA LTDA"
$newCompanyNames = @("NET-POINT", "JOSE E DA SILVA NET-POINT LTDA")
$searchFrom = 0
foreach ($newName in $newCompanyNames) {
    $rng = $d.Content
    $rng.Start = $searchFrom
    $found = $rng.Find.Execute($oldCompany, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $null, 0)
    if ($found) {
        $rng.Text = $newName
        $searchFrom = $rng.End
    }
}

# 3. CNPJ number
Replace-Unique "25.315.224/0000-11" "63.520.002/0001-21"

# 4. Registered address: the original placeholder is a lone "00" right after
#    "endereço ". Locate it precisely instead of matching "00" in isolation
#    (which would also match inside other digit strings, e.g. the CNPJ).
$anchor = $d.Content
$foundAnchor = $anchor.Find.Execute("endereço ", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, $null, 0)
if ($foundAnchor) {
    $addrRng = $d.Range($anchor.End, $anchor.End + 2)
    if ($addrRng.Text -eq "00") {
        $addrRng.Text = "R PAULO HENRIQUE MENDES, S/N, CASA 09E, TABULEIRO DO MARTINS, MACEIO - AL, CEP: 57.081-520"
    }
}

# 5. Signatory name
Replace-Unique "Alexsansdro arsaujo" "JOSE ERIVALDO DA SILVA"

# 6. Title typo fix
Replace-Unique "Sócio-Adminitrador" "Sócio-Administrador"

# 7. Email
Replace-Unique "alexsandrodro.araujo015@gmail.com" "netpoint.mcz.2025@hotmail.com"

# 8. Phone number
Replace-Unique "(88)981479415" "(82) 9609-5615"
